$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Staff_1"
$ws.Cells.Item(2, 2).Value = "DO"
$ws.Cells.Item(2, 3).Value = "M3"
$ws.Cells.Item(2, 4).Value = "M1"
$ws.Cells.Item(2, 5).Value = "M1"
$ws.Cells.Item(2, 6).Value = "M1"
$ws.Cells.Item(2, 7).Value = "A1"
$ws.Cells.Item(2, 8).Value = "M1"
$ws.Cells.Item(2, 9).Value = "DO"
$ws.Cells.Item(2, 10).Value = "M1"
$ws.Cells.Item(2, 11).Value = "M3"
$ws.Cells.Item(2, 12).Value = "M1"
$ws.Cells.Item(2, 13).Value = "M1"
$ws.Cells.Item(2, 14).Value = "A1"
$ws.Cells.Item(2, 15).Value = "M1"
$ws.Cells.Item(2, 16).Value = "DO"
$ws.Cells.Item(2, 17).Value = "M3"
$ws.Cells.Item(2, 18).Value = "M1"
$ws.Cells.Item(2, 19).Value = "M1"
$ws.Cells.Item(2, 20).Value = "M1"
$ws.Cells.Item(2, 21).Value = "A1"
$ws.Cells.Item(2, 22).Value = "M1"
$ws.Cells.Item(2, 23).Value = "DO"
$ws.Cells.Item(2, 24).Value = "M3"
$ws.Cells.Item(2, 25).Value = "M2"
$ws.Cells.Item(2, 26).Value = "M2"
$ws.Cells.Item(2, 27).Value = "M2"
$ws.Cells.Item(2, 28).Value = "A2"
$ws.Cells.Item(2, 29).Value = "M1"
$ws.Cells.Item(3, 1).Value = "Staff_2"
$ws.Cells.Item(3, 2).Value = "DO"
$ws.Cells.Item(3, 3).Value = "M1"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = "M2"
$ws.Cells.Item(3, 6).Value = "M3"
$ws.Cells.Item(3, 7).Value = "A2"
$ws.Cells.Item(3, 8).Value = "M1"
$ws.Cells.Item(3, 9).Value = "M1"
$ws.Cells.Item(3, 10).Value = "M1"
$ws.Cells.Item(3, 11).Value = "M3"
$ws.Cells.Item(3, 12).Value = "M2"
$ws.Cells.Item(3, 13).Value = "M2"
$ws.Cells.Item(3, 14).Value = "A1"
$ws.Cells.Item(3, 15).Value = "DO"
$ws.Cells.Item(3, 16).Value = "DO"
$ws.Cells.Item(3, 17).Value = "M1"
$ws.Cells.Item(3, 18).Value = "M2"
$ws.Cells.Item(3, 19).Value = "M2"
$ws.Cells.Item(3, 20).Value = "M3"
$ws.Cells.Item(3, 21).Value = "A1"
$ws.Cells.Item(3, 22).Value = "M2"
$ws.Cells.Item(3, 23).Value = "DO"
$ws.Cells.Item(3, 24).Value = "M1"
$ws.Cells.Item(3, 25).Value = "M1"
$ws.Cells.Item(3, 26).Value = "M1"
$ws.Cells.Item(3, 27).Value = "M3"
$ws.Cells.Item(3, 28).Value = "A2"
$ws.Cells.Item(3, 29).Value = "M1"
$ws.Cells.Item(4, 1).Value = "Staff_3"
$ws.Cells.Item(4, 2).Value = "DO"
$ws.Cells.Item(4, 3).Value = "M1"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = "M3"
$ws.Cells.Item(4, 6).Value = "M1"
$ws.Cells.Item(4, 7).Value = "A1"
$ws.Cells.Item(4, 8).Value = "M1"
$ws.Cells.Item(4, 9).Value = "M1"
$ws.Cells.Item(4, 10).Value = "M1"
$ws.Cells.Item(4, 11).Value = "DO"
$ws.Cells.Item(4, 12).Value = "M1"
$ws.Cells.Item(4, 13).Value = "M1"
$ws.Cells.Item(4, 14).Value = "A1"
$ws.Cells.Item(4, 15).Value = "M3"
$ws.Cells.Item(4, 16).Value = "DO"
$ws.Cells.Item(4, 17).Value = "M1"
$ws.Cells.Item(4, 18).Value = "M1"
$ws.Cells.Item(4, 19).Value = "M3"
$ws.Cells.Item(4, 20).Value = "M1"
$ws.Cells.Item(4, 21).Value = "A1"
$ws.Cells.Item(4, 22).Value = "M1"
$ws.Cells.Item(4, 23).Value = "DO"
$ws.Cells.Item(4, 24).Value = "M1"
$ws.Cells.Item(4, 25).Value = "M1"
$ws.Cells.Item(4, 26).Value = "M1"
$ws.Cells.Item(4, 27).Value = "M3"
$ws.Cells.Item(4, 28).Value = "A1"
$ws.Cells.Item(4, 29).Value = "M3"
$ws.Cells.Item(5, 1).Value = "Staff_4"
$ws.Cells.Item(5, 2).Value = "DO"
$ws.Cells.Item(5, 3).Value = "M1"
$ws.Cells.Item(5, 4).Value = "M1"
$ws.Cells.Item(5, 5).Value = "M1"
$ws.Cells.Item(5, 6).Value = "M3"
$ws.Cells.Item(5, 7).Value = "A1"
$ws.Cells.Item(5, 8).Value = "M1"
$ws.Cells.Item(5, 9).Value = "M1"
$ws.Cells.Item(5, 10).Value = "M1"
$ws.Cells.Item(5, 11).Value = "M1"
$ws.Cells.Item(5, 12).Value = "M1"
$ws.Cells.Item(5, 13).Value = "M2"
$ws.Cells.Item(5, 14).Value = "DO"
$ws.Cells.Item(5, 15).Value = "M3"
$ws.Cells.Item(5, 16).Value = "M3"
$ws.Cells.Item(5, 17).Value = "M2"
$ws.Cells.Item(5, 18).Value = "M1"
$ws.Cells.Item(5, 19).Value = "M1"
$ws.Cells.Item(5, 20).Value = "M2"
$ws.Cells.Item(5, 21).Value = "DO"
$ws.Cells.Item(5, 22).Value = "M1"
$ws.Cells.Item(5, 23).Value = "M1"
$ws.Cells.Item(5, 24).Value = "M2"
$ws.Cells.Item(5, 25).Value = "M1"
$ws.Cells.Item(5, 26).Value = "M1"
$ws.Cells.Item(5, 27).Value = "M3"
$ws.Cells.Item(5, 28).Value = "A2"
$ws.Cells.Item(5, 29).Value = "DO"
$ws.Cells.Item(6, 1).Value = "Staff_5"
$ws.Cells.Item(6, 2).Value = "DO"
$ws.Cells.Item(6, 3).Value = "M2"
$ws.Cells.Item(6, 4).Value = "M1"
$ws.Cells.Item(6, 5).Value = "M1"
$ws.Cells.Item(6, 6).Value = "M1"
$ws.Cells.Item(6, 7).Value = "A1"
$ws.Cells.Item(6, 8).Value = "M3"
$ws.Cells.Item(6, 9).Value = "DO"
$ws.Cells.Item(6, 10).Value = "M1"
$ws.Cells.Item(6, 11).Value = "M2"
$ws.Cells.Item(6, 12).Value = "M1"
$ws.Cells.Item(6, 13).Value = "M3"
$ws.Cells.Item(6, 14).Value = "A1"
$ws.Cells.Item(6, 15).Value = "M3"
$ws.Cells.Item(6, 16).Value = "DO"
$ws.Cells.Item(6, 17).Value = "M1"
$ws.Cells.Item(6, 18).Value = "M2"
$ws.Cells.Item(6, 19).Value = "M3"
$ws.Cells.Item(6, 20).Value = "M1"
$ws.Cells.Item(6, 21).Value = "A2"
$ws.Cells.Item(6, 22).Value = "M1"
$ws.Cells.Item(6, 23).Value = "DO"
$ws.Cells.Item(6, 24).Value = "M1"
$ws.Cells.Item(6, 25).Value = "M1"
$ws.Cells.Item(6, 26).Value = "M2"
$ws.Cells.Item(6, 27).Value = "M3"
$ws.Cells.Item(6, 28).Value = "A2"
$ws.Cells.Item(6, 29).Value = "M3"
$ws.Cells.Item(7, 1).Value = "Staff_6"
$ws.Cells.Item(7, 2).Value = "DO"
$ws.Cells.Item(7, 3).Value = "M1"
$ws.Cells.Item(7, 4).Value = "M1"
$ws.Cells.Item(7, 5).Value = "M3"
$ws.Cells.Item(7, 6).Value = "M1"
$ws.Cells.Item(7, 7).Value = "A1"
$ws.Cells.Item(7, 8).Value = "M1"
$ws.Cells.Item(7, 9).Value = "M1"
$ws.Cells.Item(7, 10).Value = "M1"
$ws.Cells.Item(7, 11).Value = "DO"
$ws.Cells.Item(7, 12).Value = "M1"
$ws.Cells.Item(7, 13).Value = "M3"
$ws.Cells.Item(7, 14).Value = "A1"
$ws.Cells.Item(7, 15).Value = "M1"
$ws.Cells.Item(7, 16).Value = "M1"
$ws.Cells.Item(7, 17).Value = "M1"
$ws.Cells.Item(7, 18).Value = "M1"
$ws.Cells.Item(7, 19).Value = "M3"
$ws.Cells.Item(7, 20).Value = "M1"
$ws.Cells.Item(7, 21).Value = "DO"
$ws.Cells.Item(7, 22).Value = "M1"
$ws.Cells.Item(7, 23).Value = "M1"
$ws.Cells.Item(7, 24).Value = "M1"
$ws.Cells.Item(7, 25).Value = "DO"
$ws.Cells.Item(7, 26).Value = "M1"
$ws.Cells.Item(7, 27).Value = "M1"
$ws.Cells.Item(7, 28).Value = "A1"
$ws.Cells.Item(7, 29).Value = "M3"
$ws.Cells.Item(8, 1).Value = "Staff_7"
$ws.Cells.Item(8, 2).Value = "DO"
$ws.Cells.Item(8, 3).Value = "M3"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = "M2"
$ws.Cells.Item(8, 6).Value = "M1"
$ws.Cells.Item(8, 7).Value = "A2"
$ws.Cells.Item(8, 8).Value = "M1"
$ws.Cells.Item(8, 9).Value = "M1"
$ws.Cells.Item(8, 10).Value = "M2"
$ws.Cells.Item(8, 11).Value = "M3"
$ws.Cells.Item(8, 12).Value = "M2"
$ws.Cells.Item(8, 13).Value = "M1"
$ws.Cells.Item(8, 14).Value = "DO"
$ws.Cells.Item(8, 15).Value = "M1"
$ws.Cells.Item(8, 16).Value = "M1"
$ws.Cells.Item(8, 17).Value = "DO"
$ws.Cells.Item(8, 18).Value = "M2"
$ws.Cells.Item(8, 19).Value = "M1"
$ws.Cells.Item(8, 20).Value = "M3"
$ws.Cells.Item(8, 21).Value = "A2"
$ws.Cells.Item(8, 22).Value = "M3"
$ws.Cells.Item(8, 23).Value = "M1"
$ws.Cells.Item(8, 24).Value = "M2"
$ws.Cells.Item(8, 25).Value = "M1"
$ws.Cells.Item(8, 26).Value = "M3"
$ws.Cells.Item(8, 27).Value = "DO"
$ws.Cells.Item(8, 28).Value = "A1"
$ws.Cells.Item(8, 29).Value = "M2"
$ws.Cells.Item(9, 1).Value = "Staff_8"
$ws.Cells.Item(9, 2).Value = "DO"
$ws.Cells.Item(9, 3).Value = "M1"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = "M3"
$ws.Cells.Item(9, 6).Value = "M1"
$ws.Cells.Item(9, 7).Value = "A1"
$ws.Cells.Item(9, 8).Value = "M2"
$ws.Cells.Item(9, 9).Value = "M1"
$ws.Cells.Item(9, 10).Value = "M1"
$ws.Cells.Item(9, 11).Value = "DO"
$ws.Cells.Item(9, 12).Value = "M1"
$ws.Cells.Item(9, 13).Value = "M3"
$ws.Cells.Item(9, 14).Value = "A2"
$ws.Cells.Item(9, 15).Value = "M2"
$ws.Cells.Item(9, 16).Value = "M1"
$ws.Cells.Item(9, 17).Value = "DO"
$ws.Cells.Item(9, 18).Value = "M1"
$ws.Cells.Item(9, 19).Value = "M1"
$ws.Cells.Item(9, 20).Value = "M3"
$ws.Cells.Item(9, 21).Value = "A2"
$ws.Cells.Item(9, 22).Value = "M3"
$ws.Cells.Item(9, 23).Value = "DO"
$ws.Cells.Item(9, 24).Value = "M3"
$ws.Cells.Item(9, 25).Value = "M1"
$ws.Cells.Item(9, 26).Value = "M1"
$ws.Cells.Item(9, 27).Value = "M1"
$ws.Cells.Item(9, 28).Value = "A2"
$ws.Cells.Item(9, 29).Value = "M1"
$ws.Cells.Item(10, 1).Value = "Staff_9"
$ws.Cells.Item(10, 2).Value = "M1"
$ws.Cells.Item(10, 3).Value = "M2"
$ws.Cells.Item(10, 4).Value = "M1"
$ws.Cells.Item(10, 5).Value = "M1"
$ws.Cells.Item(10, 6).Value = "M3"
$ws.Cells.Item(10, 7).Value = "DO"
$ws.Cells.Item(10, 8).Value = "M3"
$ws.Cells.Item(10, 9).Value = "DO"
$ws.Cells.Item(10, 10).Value = "M1"
$ws.Cells.Item(10, 11).Value = "M1"
$ws.Cells.Item(10, 12).Value = "M1"
$ws.Cells.Item(10, 13).Value = "M3"
$ws.Cells.Item(10, 14).Value = "A2"
$ws.Cells.Item(10, 15).Value = "M2"
$ws.Cells.Item(10, 16).Value = "DO"
$ws.Cells.Item(10, 17).Value = "M1"
$ws.Cells.Item(10, 18).Value = "M2"
$ws.Cells.Item(10, 19).Value = "M2"
$ws.Cells.Item(10, 20).Value = "M2"
$ws.Cells.Item(10, 21).Value = "A2"
$ws.Cells.Item(10, 22).Value = "M3"
$ws.Cells.Item(10, 23).Value = "DO"
$ws.Cells.Item(10, 24).Value = "M3"
$ws.Cells.Item(10, 25).Value = "M2"
$ws.Cells.Item(10, 26).Value = "M1"
$ws.Cells.Item(10, 27).Value = "M2"
$ws.Cells.Item(10, 28).Value = "A1"
$ws.Cells.Item(10, 29).Value = "M1"
